$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignUpPage")

# Append new row 7: HeadersTest / embibe.auto2@mailinator.com / embibe123 / Engineering
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "HeadersTest"

$ws.Cells.Item(7, 2).Value = "embibe.auto2@mailinator.com"
$ws.Cells.Item(7, 2).Style = "Normal"

$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "embibe123"

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "Engineering"

# Move selection to the newly added row and make this sheet the active tab
# (this also clears tabSelected on the previously-active ChooseMissionPage sheet)
$ws.Range("B7").Select()
$ws.Activate()
